# Update latest output (run 117)
# Applies changes to the "Schedule" and "Detailed" sheets of the
# optimisation_result workbook per the new run's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Schedule")
$ws2 = $wb.Worksheets.Item("Detailed")

# --- Sheet "Schedule": reshape from 2 data rows to 3 data rows ---
# Insert a new row 3 (pushes the old row 3 down to row 4), then
# rewrite every cell of rows 2-4 with the new run's values.
$ws.Rows.Item(3).Insert()

$ws.Range("A2").Value = 46042
$ws.Range("B2").Value = 46042.16666666666
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = 15.12
$ws.Range("E2").Value = 454.5780524999999
$ws.Range("F2").Value = 30.06468601190476

$ws.Range("A3").Value = 46042.33333333334
$ws.Range("B3").Value = 46042.66666666666
$ws.Range("C3").Value = 8
$ws.Range("D3").Value = 30.24
$ws.Range("E3").Value = -96.15673274999999
$ws.Range("F3").Value = -3.179786135912698

$ws.Range("E4").Value = 457.1584095
$ws.Range("F4").Value = 30.2353445436508

# --- Sheet "Detailed": update Price / Type / Pump_Status columns ---
$ws2.Range("E2").Value = "ON"
$ws2.Range("E3").Value = "ON"
$ws2.Range("E4").Value = "ON"
$ws2.Range("E5").Value = "ON"
$ws2.Range("E6").Value = "ON"
$ws2.Range("B7").Value = 57.06003
$ws2.Range("E7").Value = "ON"
$ws2.Range("B8").Value = 57.06003
$ws2.Range("E8").Value = "ON"
$ws2.Range("B9").Value = 64.89
$ws2.Range("C9").Value = "historical"
$ws2.Range("E9").Value = "ON"
$ws2.Range("B10").Value = 65.8814
$ws2.Range("C10").Value = "historical"
$ws2.Range("E10").Value = "OFF"
$ws2.Range("B11").Value = 67.01821
$ws2.Range("E11").Value = "OFF"
$ws2.Range("B12").Value = 77.94
$ws2.Range("E12").Value = "OFF"
$ws2.Range("B13").Value = 81.61441
$ws2.Range("E13").Value = "OFF"
$ws2.Range("B14").Value = 77.94
$ws2.Range("E14").Value = "OFF"
$ws2.Range("B15").Value = 77.94
$ws2.Range("E15").Value = "OFF"
$ws2.Range("E16").Value = "OFF"
$ws2.Range("B17").Value = 8.66036
$ws2.Range("E17").Value = "OFF"
$ws2.Range("B18").Value = 0.7
$ws2.Range("B19").Value = 0.0099
$ws2.Range("B20").Value = -1.17663
$ws2.Range("B21").Value = -5.50985
$ws2.Range("B22").Value = -6.57149
$ws2.Range("B23").Value = -6.78314
$ws2.Range("B24").Value = -8.26613
$ws2.Range("B25").Value = -7.64915
$ws2.Range("B26").Value = -6.0195
$ws2.Range("B27").Value = -8.86058
$ws2.Range("B28").Value = -9.83442
$ws2.Range("B29").Value = -7.95321
$ws2.Range("B30").Value = -7.73429
$ws2.Range("B31").Value = -9.0256
$ws2.Range("B32").Value = -7.49011
$ws2.Range("B33").Value = -6.45809
$ws2.Range("B34").Value = -4.87045
$ws2.Range("B36").Value = -5.08833
$ws2.Range("B37").Value = 5.09213
$ws2.Range("B38").Value = 9.64796
$ws2.Range("B39").Value = 33.25863
$ws2.Range("B40").Value = 55.50567
$ws2.Range("B41").Value = 57.20458
$ws2.Range("B42").Value = 61.81708
$ws2.Range("B43").Value = 57.88272
$ws2.Range("B44").Value = 57.3
$ws2.Range("B45").Value = 57.76776
$ws2.Range("B46").Value = 57.06007
$ws2.Range("B47").Value = 58.64616
$ws2.Range("B48").Value = 61.10663
$ws2.Range("B49").Value = 57.3
